# Disaggregated ISIC spending allocations (except construction)
# On the "Pre ISIC Consolidation" sheet, row 7 ("EU ISIC Groupings") held several
# cells that were all lumped into a single (coarser) ISIC code. Break these out
# into their own, more granular ISIC codes (mirroring the labels already used
# in row 2 "Default EPS ISIC Groupings").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre ISIC Consolidation")

$ws.Range("E7").Value  = "ISIC 07T08"
$ws.Range("N7").Value  = "ISIC 22"
$ws.Range("S7").Value  = "ISIC 25"
$ws.Range("T7").Value  = "ISIC 26"
$ws.Range("U7").Value  = "ISIC 27"
$ws.Range("X7").Value  = "ISIC 30"

$wb.Application.Calculate()
